# Esquema de barras modificado
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New header cells on row 8 (Z8:AD8) ---
$ws.Range("Z8").Value = "prx"
$ws.Range("AA8").Value = "alpha"
$ws.Range("AB8").Value = "r"
$ws.Range("AC8").Value = "y"
$ws.Range("AD8").Value = "z"

# --- Row 9 ---
$ws.Range("A9").Value = 2
$ws.Range("M9").Formula = '=X7-H10-$X$16'
$ws.Range("O9").Formula = '=-X8+$X$16'

$ws.Range("Z9").Value = 1
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Formula = "=X13/2 - X13/4"
$ws.Range("AC9").Formula = "=AB9*COS(AA9)"
$ws.Range("AD9").Formula = "=AB9*SIN(AA9)"
$ws.Range("AC9:AD9").NumberFormat = "0.00000"

# --- Row 10 ---
$ws.Range("L10").Value = 0
$ws.Range("O10").Formula = '=X10-$X$16'

$ws.Range("Z10").Value = 2
$ws.Range("AA10").Formula = "=AA9+120"
$ws.Range("AB10").Formula = "=AB9"
$ws.Range("AC10").Formula = "=-AB10*COS(AA10)"
$ws.Range("AD10").Formula = "=AB10*SIN(AA10)"
$ws.Range("AC10:AD10").NumberFormat = "0.00000"

# --- Row 11 ---
$ws.Range("Z11").Value = 3
$ws.Range("AA11").Formula = "=AA10+120"
$ws.Range("AB11").Formula = "=AB10"
$ws.Range("AC11").Formula = "=-AB11*COS(AA11)"
$ws.Range("AD11").Formula = "=-AB11*SIN(AA11)"
$ws.Range("AC11:AD11").NumberFormat = "0.00000"

# --- Row 13 ---
$ws.Range("B13").Formula = "=X12/2+X14/2"

# --- View state ---
$ws.Range("AF5").Select()
